$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.711.94"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.507.58"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'322.56"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'108.02"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.559"
$ws.Range("E9").Value = "  +2.90%  "
$ws.Range("D10").Value = "'40.37"
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").Value = "'19.52"
$ws.Range("E12").Value = "  +5.02%  "
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").Value = "'7.16"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "2.900.53"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "2.515.30"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "'0.850"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "47.625.20"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "'13.31"
$ws.Range("E19").Value = "  +2.10%  "
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("B21").Value = "ImmutableX"
$ws.Range("C21").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D21").Value = "'2.78"
$ws.Range("E21").Value = "  +8.93%  "
$ws.Range("B22").Value = "ShibaInu"
$ws.Range("C22").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D22").Value = "0.0₃0940"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").Value = "'70.96"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'247.43"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'25.76"
$ws.Range("E27").Value = "  -1.69%  "
$ws.Range("D28").Value = "'10.19"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.142"
$ws.Range("E29").Value = "  +3.85%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "'34.87"
$ws.Range("E30").Value = "  -2.72%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.08"
$ws.Range("E31").Value = "  -9.37%  "
$ws.Range("D32").Value = "'49.85"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "'20.01"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").Value = "'5.38"
$ws.Range("E34").Value = "  -1.49%  "
$ws.Range("D36").Value = "'0.0784"
$ws.Range("E36").Value = "  -1.20%  "
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "'22.41"
$ws.Range("E41").Value = "  +4.68%  "
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").Value = "'118.71"
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("D45").Value = "2.006.23"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D46").Value = "'3.10"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").Value = "'1.80"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").Value = "'5.16"
$ws.Range("E50").Value = "  -3.72%  "
$ws.Range("D51").Value = "'56.50"
$ws.Range("E51").Value = "  +1.82%  "
